$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H1").Value = "length_pipe"
$ws.Range("A3").Value = "length_pipe"

$ws.Range("B2").Value = 0.0001315260047151133
$ws.Range("C2").Value = 1.299046696460312
$ws.Range("D2").Value = 0.8534736798886241
$ws.Range("E2").Value = 1.472183622241831
$ws.Range("F2").Value = -12.24395991558717
$ws.Range("G2").Value = 3.468721592776412
$ws.Range("K2").Value = 0.0196
$ws.Range("P2").Value = -0.3917021426743598
$ws.Range("Q2").Value = -0.1039204065152347

$ws.Range("B3").Value = 0.0000006129859035316447
$ws.Range("G3").Value = 3.468721592776412
$ws.Range("H3").Value = 16.95283206967629
$ws.Range("K3").Value = 0.0196

$ws.Range("B4").Value = 0.0000005610385877157079
$ws.Range("G4").Value = 3.174724336146332
$ws.Range("K4").Value = 0.0196

$ws.Range("B5").Value = 0.000001220899222172234
$ws.Range("G5").Value = 6.90979573
$ws.Range("J5").Value = 7.281418596835137
$ws.Range("K5").Value = 0.0196

$ws.Range("B6").Value = 0.0000007091473238619253
$ws.Range("G6").Value = 3.468721592776412
$ws.Range("K6").Value = 0.02267528804567081

$ws.Range("B7").Value = 0.0000003846624453452384
$ws.Range("G7").Value = 3.468721592776412
$ws.Range("K7").Value = 0.0196

$ws.Range("B8").Value = 0.000000961911140693035
$ws.Range("G8").Value = 3.468721592776412
$ws.Range("K8").Value = 0.0196

$ws.Range("B9").Value = 0.000001262852011812774
$ws.Range("G9").Value = 3.468721592776412
$ws.Range("K9").Value = 0.0196

$ws.Range("B10").Value = 0.0000005612328507227877
$ws.Range("G10").Value = 3.468721592776412
$ws.Range("K10").Value = 0.0196

$ws.Range("B11").Value = 0.0000004857380916276359
$ws.Range("G11").Value = 3.468721592776412
$ws.Range("K11").Value = 0.0196

$ws.Range("B12").Value = 0.0000004933410007644325
$ws.Range("G12").Value = 3.468721592776412
$ws.Range("K12").Value = 0.0196

$ws.Range("B13").Value = 0.0000004805405526785011
$ws.Range("G13").Value = 3.468721592776412
$ws.Range("K13").Value = 0.0196
